# Weekly update: insert two new daily price records (rows 357-358) for
# "Pepino ensalada" at Vega Central Mapocho de Santiago, pushing the
# previously existing rows 357-377 down to 359-379.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 357 and 358, shifting rows 357:377 down to 359:379.
$ws.Range("A357:A358").EntireRow.Insert()

# --- New row 357 ---
$ws.Cells.Item(357, 1).Value  = 9
$ws.Cells.Item(357, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(357, 3).Value  = "Metropolitana"
$ws.Cells.Item(357, 4).Value  = 44931
$ws.Cells.Item(357, 5).Value  = 13
$ws.Cells.Item(357, 6).Value  = 100112043
$ws.Cells.Item(357, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(357, 8).Value  = "Sin especificar"
$ws.Cells.Item(357, 9).Value  = "Primera"
$ws.Cells.Item(357, 10).Value = 70
$ws.Cells.Item(357, 11).Value = 16000
$ws.Cells.Item(357, 12).Value = 17000
$ws.Cells.Item(357, 13).Value = 16500
$ws.Cells.Item(357, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(357, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(357, 16).Value = 330
$ws.Cells.Item(357, 17).Value = 50
$ws.Cells.Item(357, 18).Value = "Hortaliza"

# --- New row 358 ---
$ws.Cells.Item(358, 1).Value  = 9
$ws.Cells.Item(358, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(358, 3).Value  = "Metropolitana"
$ws.Cells.Item(358, 4).Value  = 44931
$ws.Cells.Item(358, 5).Value  = 13
$ws.Cells.Item(358, 6).Value  = 100112043
$ws.Cells.Item(358, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(358, 8).Value  = "Sin especificar"
$ws.Cells.Item(358, 9).Value  = "Primera"
$ws.Cells.Item(358, 10).Value = 90
$ws.Cells.Item(358, 11).Value = 13000
$ws.Cells.Item(358, 12).Value = 15000
$ws.Cells.Item(358, 13).Value = 14000
$ws.Cells.Item(358, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(358, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(358, 16).Value = 233
$ws.Cells.Item(358, 17).Value = 60
$ws.Cells.Item(358, 18).Value = "Hortaliza"
